$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Fill in laser/filter information (columns C and D) for rows that were
# missing it, following the pattern already used by rows 62-76, for lot AJ02.

# Row 77: BV421
$ws.Range("C77").Value = 405
$ws.Range("D77").Value = "450/50"

# Row 79: BV605
$ws.Range("C79").Value = 405
$ws.Range("D79").Value = "610/20"

# Row 83: FITC
$ws.Range("C83").Value = 488
$ws.Range("D83").Value = "530/40"

# Row 84: PerCP-Cy5.5
$ws.Range("C84").Value = 560
$ws.Range("D84").Value = "695/40"

# Row 85: PE
$ws.Range("C85").Value = 488
$ws.Range("D85").Value = "575/25"

# Row 87: PE-Cy5 excited at 488 -> rename label, add filter info
$ws.Range("B87").Value = "PE-Cy5 488"
$ws.Range("D87").Value = "680/30"

# Row 88: PE-Cy5 excited at 561 -> rename label
$ws.Range("B88").Value = "PE-Cy5 561"

# Row 90: PE-Cy7
$ws.Range("C90").Value = 488
$ws.Range("D90").Value = "750LP"

# Row 91: APC
$ws.Range("C91").Value = 633
$ws.Range("D91").Value = "665/20"

# Row 93: APC-Cy7
$ws.Range("C93").Value = 635
$ws.Range("D93").Value = "750LP"

# Update the view/selection state to reflect where the user was working
$excel.ActiveWindow.ScrollRow = 63
[void]$ws.Range("C93:D93").Select()
